$d = $word.ActiveDocument

function Set-ParaText($searchText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "text not found: $searchText"
    }
    $rng.Text = $newText
    return $rng
}

function Insert-ItalicParaAfter($rng, $newText) {
    # rng must currently cover exactly the paragraph whose text was just set
    $paraIndex = $rng.Paragraphs(1).Index
    $rng.InsertParagraphAfter() | Out-Null
    $newRng = $d.Paragraphs($paraIndex + 1).Range
    $newRng.Text = $newText
    # Exclude the paragraph mark from the italic formatting so the new
    # paragraph doesn't pick up an <w:pPr><w:rPr><w:i/></w:rPr></w:pPr>.
    $newRng.MoveEnd(1, -1) | Out-Null
    $newRng.Font.Italic = 1
}

# 1) Ativação date
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2) | Out-Null

# 2) Objetivos paragraph + new italic English paragraph
$r = Set-ParaText "Fornecer uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base, de transformação." "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."
Insert-ItalicParaAfter $r "Provide students with a current view of industrial processes that use chemical conversion as a route to transform raw material into product. The processes of the basic chemical and transformation industries will be studied."

# 3) Programa resumido paragraph + new italic English paragraph
$r = Set-ParaText "Introdução ao Estudo dos Processos Químicos Industriais. Relacionamento com a Engenharia Química; Derivados Inorgânicos do Nitrogênio; Ácido Sulfúrico; Fósforo e Ácido Fosfórico; Fertilizantes; Indústrias de Cloro Álcalis; Indústrias de Vidros e Materiais Cerâmicos; Carga e Pigmentos Inorgânicos." "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"
Insert-ItalicParaAfter $r "Introduction to Industrial Chemical Processes; NPK / Fertilizers; Sulfuric Acid; Chlorine Alkali; Paper and Cellulose; Sugar and alcohol; Biotechnological Processes."

# 4) Programa paragraph + new italic English paragraph
$r = Set-ParaText "Introdução ao estudo dos Processos Químicos Industriais. Relacionamento com a Engenharia Química. Fundamentos dos processos químicos. Condução dos processos (batelada X contínuo). Fluxogramas. Derivados inorgânicos do nitrogênio - Introdução Amônia. Generalidades. Amônia. Produção sintética pelo processo Haber Bosch. Uréia: Generalidades. Processo de Fabricação. Nitrato de Amônia: Generalidades - Processo de Fabricação. Acído Nítrico. Generalidades. Processo de Fabricação Ácido Sulfúrico. Generalidades. Processo de Fabricação. Concentração. Fósforo e Ácido fosfórico. Generalidades. Matérias Primas. Produção de ácido fosfórico.  Indústrias de cloro álcalis. Generalidades. Matérias primas. Produção de barrilha e bicarbonato de sódio. Indústria de cloro e álcalis: produção de cloro e soda caústica. Células a diafragma.e mercúrio. Comparação. Ácido Clorídrico: fabricação e aplicações. Indústrias de vidros e materiais cerâmicos: matérias primas e fabricação." ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."
Insert-ItalicParaAfter $r "1. Introduction to Industrial Chemical Processes;2. NPK / Fertilizers3. Sulfuric Acid;4. Chlorine Alkali;5. Paper and Cellulose;6. Sugar and alcohol;7. Biotechnological Processes;"

# 5) Método value
$d.Content.Find.Execute("Serão aplicadas duas provas (P1 e P2) e a Nota Final (NF) será a média aritmética das mesmas.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos.", 2) | Out-Null

# 6) Critério value
$d.Content.Find.Execute("Serão aprovados os alunos com NF maior ou igual a 5,0 e frequência superior a 70%.", $true, $false, $false, $false, $false, $true, 1, $false, "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula.", 2) | Out-Null

# 7) Norma de recuperação value
Set-ParaText "Será feita a Recuperação( REC) para alunos com NF maior ou igual a 3,0 e menor que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou superior a 5,0, sendo MF= (NF+ REC)/2 .                                                 Na semana anterior à REC será dada uma aula de recordação de toda a matéria apresentada." "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação." | Out-Null

# 8) Bibliografia paragraph - whole paragraph (multiple runs/breaks) collapses to one run
$bibHeaderRng = $d.Content
$bibHeaderRng.Find.Execute("Bibliografia") | Out-Null
$bibParaIndex = $bibHeaderRng.Paragraphs(1).Index
$bibRng = $d.Paragraphs($bibParaIndex + 1).Range
$bibRng.Text = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"

Write-Host "edit complete"
